$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DataRow {
    param($RowIndex, $Values)

    $ws.Cells.Item($RowIndex, 1).Value  = $Values[0]   # Mercado ID
    $ws.Cells.Item($RowIndex, 2).Value  = $Values[1]   # Mercado
    $ws.Cells.Item($RowIndex, 3).Value  = $Values[2]   # Región
    $ws.Cells.Item($RowIndex, 4).Value  = $Values[3]   # Fecha (serial)
    $ws.Cells.Item($RowIndex, 5).Value  = $Values[4]   # Codreg
    $ws.Cells.Item($RowIndex, 6).Value  = $Values[5]   # Categoría ID
    $ws.Cells.Item($RowIndex, 7).Value  = $Values[6]   # Categoría
    $ws.Cells.Item($RowIndex, 8).Value  = $Values[7]   # Variedad
    $ws.Cells.Item($RowIndex, 9).Value  = $Values[8]   # Calidad
    $ws.Cells.Item($RowIndex, 10).Value = $Values[9]   # Volumen
    $ws.Cells.Item($RowIndex, 11).Value = $Values[10]  # Precio mínimo
    $ws.Cells.Item($RowIndex, 12).Value = $Values[11]  # Precio máximo
    $ws.Cells.Item($RowIndex, 13).Value = $Values[12]  # Precio promedio ponderado
    $ws.Cells.Item($RowIndex, 14).Value = $Values[13]  # Unidad de comercialización
    $ws.Cells.Item($RowIndex, 15).Value = $Values[14]  # Origen
    $ws.Cells.Item($RowIndex, 16).Value = $Values[15]  # Precio $/Kg
    $ws.Cells.Item($RowIndex, 17).Value = $Values[16]  # Kg o Unidades
    $ws.Cells.Item($RowIndex, 18).Value = $Values[17]  # Clasificación
}

# --- Insert new row at 6 (pushes old row 6 onward down by one) ---
$ws.Rows.Item(6).Insert()
Set-DataRow 6 @(8, "Terminal La Palmera de La Serena", "Coquimbo", 44847, 4, 100112026, "Haba", "Sin especificar", "Primera", 520, 7000, 8000, 7500, "`$/saco 25 kilos", "Provincia del Elquí", 300, 25, "Hortaliza")

# --- Insert new row at 15 (after the first insertion has shifted everything) ---
$ws.Rows.Item(15).Insert()
Set-DataRow 15 @(8, "Terminal La Palmera de La Serena", "Coquimbo", 44848, 4, 100112026, "Haba", "Sin especificar", "Primera", 800, 7000, 8000, 7500, "`$/saco 25 kilos", "Provincia del Elquí", 300, 25, "Hortaliza")
